$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.941.74"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.524.03"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.32"
$ws.Range("E5").Value = "  +1.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.21"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.523.90"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.86"
$ws.Range("E11").Value = "  -3.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  +2.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.112.85"
$ws.Range("E13").Value = "  +2.43%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.14"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.523.38"
$ws.Range("E16").Value = "  +2.44%  "

$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.888.40"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.03"
$ws.Range("E19").Value = "  +3.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.82"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.17"
$ws.Range("E21").Value = "  +5.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.97"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.571"
$ws.Range("E23").Value = "  +3.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.661.09"
$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.65"
$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +4.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  +6.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").Value = "  +1.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.533.62"
$ws.Range("E32").Value = "  +2.64%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  +16.57%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.75"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.143"
$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "168.96"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  +6.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.84"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.94"
$ws.Range("E40").Value = "  +5.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  +5.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.818"
$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.14"
$ws.Range("E43").Value = "  +14.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.30"
$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.41"
$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("E47").Value = "  +5.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.66"
$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.83"
$ws.Range("E49").Value = "  +5.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.388.88"
$ws.Range("E50").Value = "  +10.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "305.00"
$ws.Range("E51").Value = "  +13.34%  "
